# The NATMI pipeline was rerun with updated per-cell TPM values.
# Only the Egf (ligand) / Erbb2 (receptor) expression figures for the
# "ECs" cluster changed; this shifts the detection rates, specificity
# scores and edge weights recomputed from them across every data row
# (rows 2-10, columns E:T) of the Egf-Erbb2 LR-pair sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values, one row per data record, columns E..T in order.
$newValues = @{
    2 = @(3, 1, 0.04983966666666667, 0.149519, 0.1823731600337622, 0.1823731600337622, 3, 1, 2.451416666666667, 7.35425, 0.2191928499183569, 0.2191928499183569, 0.1221777895277778, 1.09960010575, 0.03997489269641692, 0.03997489269641692)
    3 = @(3, 1, 0.04983966666666667, 0.149519, 0.1823731600337622, 0.1823731600337622, 3, 1, 4.973328333333334, 14.919985, 0.4446889938320204, 0.4446889938320204, 0.2478690263572223, 2.230821237215, 0.08109933703737976, 0.08109933703737975)
    4 = @(3, 1, 0.04983966666666667, 0.149519, 0.1823731600337622, 0.1823731600337622, 3, 1, 3.75909, 11.27727, 0.3361181562496228, 0.3361181562496228, 0.18735179257, 1.68616613313, 0.06129893029996555, 0.06129893029996554)
    5 = @(3, 1, 0.1101723333333333, 0.3305169999999999, 0.4031422744592926, 0.4031422744592926, 3, 1, 2.451416666666667, 7.35425, 0.2191928499183569, 0.2191928499183569, 0.2700782941388888, 2.43070464725, 0.08836590406130077, 0.08836590406130075)
    6 = @(3, 1, 0.1101723333333333, 0.3305169999999999, 0.4031422744592926, 0.4031422744592926, 3, 1, 4.973328333333334, 14.919985, 0.4446889938320204, 0.4446889938320204, 0.5479231869161111, 4.931308682245, 0.1792729324004551, 0.179272932400455)
    7 = @(3, 1, 0.1101723333333333, 0.3305169999999999, 0.4031422744592926, 0.4031422744592926, 3, 1, 3.75909, 11.27727, 0.3361181562496228, 0.3361181562496228, 0.4141477165099999, 3.727329448589999, 0.1355034379975368, 0.1355034379975368)
    8 = @(3, 1, 0.113272, 0.339816, 0.4144845655069452, 0.4144845655069451, 3, 1, 2.451416666666667, 7.35425, 0.2191928499183569, 0.2191928499183569, 0.2776768686666667, 2.499091818, 0.0908520531606392, 0.09085205316063918)
    9 = @(3, 1, 0.113272, 0.339816, 0.4144845655069452, 0.4144845655069451, 3, 1, 4.973328333333334, 14.919985, 0.4446889938320204, 0.4446889938320204, 0.5633388469733334, 5.07004962276, 0.1843167243941856, 0.1843167243941856)
    10 = @(3, 1, 0.113272, 0.339816, 0.4144845655069452, 0.4144845655069451, 3, 1, 3.75909, 11.27727, 0.3361181562496228, 0.3361181562496228, 0.42579964248, 3.83219678232, 0.1393157879521204, 0.1393157879521204)
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $ws.Cells.Item($row, 5 + $i).Value = $rowValues[$i]
    }
}
